$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting rows 7..34 down to 8..35
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new data point
$ws.Cells.Item(7, 1).Value = 8
$ws.Cells.Item(7, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 100112052
$ws.Cells.Item(7, 7).Value = "Albahaca"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 800
$ws.Cells.Item(7, 11).Value = 4500
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 13).Value = 4750
$ws.Cells.Item(7, 14).Value = '$/paquete'
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 4750
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# Copy the date number format from the row above so the new row matches existing formatting
$ws.Cells.Item(8, 4).Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4122)
